$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (Type) to hold the "Variable" column
$ws.Columns("B:B").Insert()

# Header for new column
$ws.Range("B1").Value = "Variable"

# Values for new column (variable name "c1" for each translatable row)
$ws.Range("B2").Value = "c1"
$ws.Range("B3").Value = "c1"
$ws.Range("B4").Value = "c1"
$ws.Range("B5").Value = "c1"

$ws.Columns("A:F").AutoFit()

# Match the resulting selection state (Excel leaves selection at B6 after edits)
[void]$ws.Range("B6").Select()
